$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (index 3), shifting existing columns C:F to D:G
$ws.Columns.Item(3).Insert()

# Match the new column's width to column B (its left neighbour), mirroring
# Excel's default "insert column" formatting behaviour
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Set the header for the newly inserted column
$ws.Range("C1").Value = "Email"

# Update selection to reflect where the edit was made
$ws.Range("C1").Select()
